$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-23 Thursday", "2025-01-24 Friday"),
    @("228÷9=25, 3", "419÷2=209, 1"),
    @("158÷7=22, 4", "662÷9=73, 5"),
    @("831÷9=92, 3", "188÷2=94, 0"),
    @("334÷7=47, 5", "984÷7=140, 4"),
    @("739÷4=184, 3", "710÷5=142, 0"),
    @("970÷6=161, 4", "923÷3=307, 2"),
    @("182÷6=30, 2", "706÷6=117, 4"),
    @("535÷5=107, 0", "892÷9=99, 1"),
    @("315÷4=78, 3", "482÷5=96, 2"),
    @("870÷8=108, 6", "549÷3=183, 0"),
    @("150÷3=50, 0", "306÷7=43, 5"),
    @("223÷5=44, 3", "448÷7=64, 0"),
    @("619÷5=123, 4", "401÷9=44, 5"),
    @("755÷4=188, 3", "421÷2=210, 1"),
    @("423÷9=47, 0", "427÷3=142, 1"),
    @("902÷9=100, 2", "209÷9=23, 2"),
    @("496÷5=99, 1", "207÷4=51, 3"),
    @("671÷9=74, 5", "181÷6=30, 1"),
    @("290÷8=36, 2", "673÷2=336, 1"),
    @("673÷8=84, 1", "425÷8=53, 1"),
    @("271÷9=30, 1", "780÷5=156, 0"),
    @("564÷2=282, 0", "325÷6=54, 1"),
    @("918÷5=183, 3", "443÷9=49, 2"),
    @("324÷9=36, 0", "653÷7=93, 2"),
    @("895÷3=298, 1", "811÷2=405, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
